$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Coin/Link/Price/Volume columns for the updated rows keep their original
# text representation (e.g. "107.60", "0.120") instead of Excel auto-coercing
# numeric-looking strings into trimmed numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.928.00"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +6.95%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.248.58"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.89%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "394.09"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.65%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.60"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.66%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.246.91"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.74%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.566"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.24%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.02%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.617"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.24%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "38.94"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.62%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0971"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +11.93%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.87%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.745.91"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.65%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.14"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.84%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.94"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.07%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.244.56"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.99%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.03"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.78%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.74"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.12%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "56.648.84"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +6.51%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.30"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.22%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000105"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +8.17%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.95"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.46%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "298.76"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +10.09%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.56"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.28%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.13"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.45%  "

# Row 27
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.88"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.78%  "

# Row 28
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.38"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.90%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.74"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.46%  "

# Row 30
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.20"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.69%  "

# Row 31
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.168"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.10%  "

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.03%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.11%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.95"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.12%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "37.21"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.98%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0484"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.06%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.11"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.40%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "51.26"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.46%  "

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.53"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.29%  "

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.52%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.05"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +11.11%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "134.10"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.31%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.90"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.46%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.120"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.71%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.98"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.02%  "

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.284"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.91%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.92"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.79%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.87"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.61%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.138.20"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.41%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.47%  "

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.04"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +20.80%  "

